# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (and before the
#    "总计" summary sheet), populated with the quarter's fund-holding row -
#    same layout/formatting as the other per-quarter sheets (e.g. "2021-Q4").
# 2. Update the "总计" (summary) sheet: add a new top data row for 2022-Q1
#    and shift the existing rows down, renumbering the index column.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Step 1: new "2022-Q1" sheet, inserted before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q4Sheet = $wb.Worksheets.Item("2021-Q4")

$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# NOTE: the sheet object passed as the "before" argument to Worksheets.Add()
# is left stale afterwards (it ends up pointing at the newly-inserted sheet
# instead of the original one) - re-resolve "总计" by name before using it
# again.
$totalSheet = $wb.Worksheets.Item("总计")

# Seed formatting by copying the "2021-Q4" sheet's header + data row, then
# overwrite with this quarter's values.
$q4Sheet.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial($xlPasteFormats)
$q4Sheet.Range("A2:H2").Copy()
$q1.Range("A2:H2").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0

$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "513080"
$q1.Range("B2").Style = "Normal"

$q1.Range("C2").Value = "华安法国CAC40ETF（QDII）"
$q1.Range("C2").Style = "Normal"

$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "0.60"
$q1.Range("D2").Style = "Normal"

$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "96.69"
$q1.Range("E2").Style = "Normal"

$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "6.86"
$q1.Range("F2").Style = "Normal"

$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.0412"
$q1.Range("G2").Style = "Normal"

$q1.Range("H2").Value = 3

# ---------------------------------------------------------------------
# Step 2: update "总计" sheet with the new quarter + shifted rows
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$dates = @("2022-Q1", "2021-Q4", "2021-Q3", "2021-Q2", "2021-Q1", "2020-Q4")
$counts = @(1, 1, 1, 1, 1, 2)
$values = @(0.04, 0.04, 0.03, 0.04, 0.03, 0.04)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $i + 2
    $totalSheet.Cells.Item($r, 1).Value = $i
    $totalSheet.Cells.Item($r, 2).Value = $dates[$i]
    $totalSheet.Cells.Item($r, 3).Value = $counts[$i]
    $totalSheet.Cells.Item($r, 4).Value = $values[$i]
}

# The row Insert() copied the header's bold/bordered formatting onto the new
# row; restore the plain data-row look (column A keeps the bordered style).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial($xlPasteFormats)
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0
